$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01675466666666667
$ws.Range("H2").Value = 0.050264
$ws.Range("I2").Value = 0.0001854906931657378
$ws.Range("J2").Value = 0.0001854906931657378
$ws.Range("M2").Value = 8.142376000000001
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 0.1364227957546667
$ws.Range("R2").Value = 1.227805161792
$ws.Range("S2").Value = 0.00003229975285025907
$ws.Range("T2").Value = 0.00003229975285025907
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01675466666666667
$ws.Range("H3").Value = 0.050264
$ws.Range("I3").Value = 0.0001854906931657378
$ws.Range("J3").Value = 0.0001854906931657378
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 0.4078143558568889
$ws.Range("R3").Value = 3.670329202712
$ws.Range("S3").Value = 0.00009655499896551947
$ws.Range("T3").Value = 0.00009655499896551946
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01675466666666667
$ws.Range("H4").Value = 0.050264
$ws.Range("I4").Value = 0.0001854906931657378
$ws.Range("J4").Value = 0.0001854906931657378
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 0.239210296592
$ws.Range("R4").Value = 2.152892669328
$ws.Range("S4").Value = 0.00005663594134995925
$ws.Range("T4").Value = 0.00005663594134995924
$ws.Range("I5").Value = 0.9933938536206305
$ws.Range("J5").Value = 0.9933938536206304
$ws.Range("M5").Value = 8.142376000000001
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 730.611140017352
$ws.Range("R5").Value = 6575.500260156168
$ws.Range("S5").Value = 0.1729810558540708
$ws.Range("T5").Value = 0.1729810558540708
$ws.Range("I6").Value = 0.9933938536206305
$ws.Range("J6").Value = 0.9933938536206304
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("S6").Value = 0.5170994882368057
$ws.Range("T6").Value = 0.5170994882368057
$ws.Range("I7").Value = 0.9933938536206305
$ws.Range("J7").Value = 0.9933938536206304
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("S7").Value = 0.3033133095297539
$ws.Range("T7").Value = 0.3033133095297538
$ws.Range("G8").Value = 0.5799533333333334
$ws.Range("I8").Value = 0.006420655686203657
$ws.Range("J8").Value = 0.006420655686203655
$ws.Range("M8").Value = 8.142376000000001
$ws.Range("N8").Value = 24.427128
$ws.Range("O8").Value = 0.1741313933276368
$ws.Range("P8").Value = 0.1741313933276368
$ws.Range("Q8").Value = 4.722198102453334
$ws.Range("R8").Value = 42.49978292208001
$ws.Range("S8").Value = 0.001118037720715657
$ws.Range("T8").Value = 0.001118037720715656
$ws.Range("G9").Value = 0.5799533333333334
$ws.Range("I9").Value = 0.006420655686203657
$ws.Range("J9").Value = 0.006420655686203655
$ws.Range("O9").Value = 0.5205382400466131
$ws.Range("P9").Value = 0.5205382400466131
$ws.Range("S9").Value = 0.003342196810841731
$ws.Range("T9").Value = 0.00334219681084173
$ws.Range("G10").Value = 0.5799533333333334
$ws.Range("I10").Value = 0.006420655686203657
$ws.Range("J10").Value = 0.006420655686203655
$ws.Range("O10").Value = 0.3053303666257501
$ws.Range("P10").Value = 0.3053303666257501
$ws.Range("Q10").Value = 8.280129449080002
$ws.Range("R10").Value = 74.52116504172001
$ws.Range("S10").Value = 0.00196042115464627
$ws.Range("T10").Value = 0.001960421154646269
